$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3.0TD COMP")

# Update the "precio de potencia" coefficient (J2:J7) for every tariff
# period on the "3.0TD COMP" sheet: 0.0442 -> 0.04325. Every downstream
# formula (L, M, N columns, the summary block in rows 11-23, etc.) is
# recalculated automatically by the engine.
$ws.Range("J2:J7").Value = 0.04325

# Leave the selection where the author left it after editing the range.
$ws.Range("J2:J7").Select() | Out-Null
